# Update computed CG values on "GLOBAL RESULTS", "FUEL TANK" and
# "LANDING GEARS" sheets following a re-run of the balance computation
# ("Debugging and Add tesi Spoti").

$wb = $excel.ActiveWorkbook

# --- GLOBAL RESULTS -------------------------------------------------
$ws = $wb.Worksheets.Item("GLOBAL RESULTS")

$ws.Range("C3").Value  = 17.246254408436393
$ws.Range("C5").Value  = -0.7344285755648581
$ws.Range("C7").Value  = 34.24872520773436
$ws.Range("C9").Value  = -18.83989776057891
$ws.Range("C13").Value = 16.411636064310827
$ws.Range("C15").Value = -0.9325556693060377
$ws.Range("C17").Value = 12.838713675522692
$ws.Range("C19").Value = -23.922344595948303
$ws.Range("C23").Value = 16.411636064310827
$ws.Range("C25").Value = -0.9325556693060377
$ws.Range("C27").Value = 12.838713675522692
$ws.Range("C29").Value = -23.922344595948303
$ws.Range("C33").Value = 16.411636064310827
$ws.Range("C35").Value = -0.9325556693060377
$ws.Range("C37").Value = 12.838713675522692
$ws.Range("C39").Value = -23.922344595948303
$ws.Range("C43").Value = 16.92960575115403
$ws.Range("C45").Value = -0.6896245039774538
$ws.Range("C47").Value = 26.12590901061653
$ws.Range("C49").Value = -17.690563222070327
$ws.Range("C53").Value = 16.82667937481264
$ws.Range("C55").Value = -0.8455920374946073
$ws.Range("C57").Value = 23.485594524422805
$ws.Range("C59").Value = -21.69151373406922
$ws.Range("C62").Value = 7.0208304745943515
$ws.Range("C63").Value = 26.12590901061653
$ws.Range("C64").Value = 29.619176641734317

# --- FUEL TANK --------------------------------------------------------
$ws = $wb.Worksheets.Item("FUEL TANK")

$ws.Range("C2").Value = 4.389801741976934
$ws.Range("C6").Value = 16.559801741976933

# --- LANDING GEARS ------------------------------------------------------
$ws = $wb.Worksheets.Item("LANDING GEARS")

$ws.Range("C2").Value  = 16.920581453544408
$ws.Range("C6").Value  = 16.920581453544404
$ws.Range("C11").Value = 16.920581453544408
